$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.43217854801802
$ws.Range("C2").Value = 11.0177613601999
$ws.Range("D2").Value = 5.058263248463851
$ws.Range("F2").Value = 23.68727751065412
$ws.Range("G2").Value = 27.68821453509152
$ws.Range("H2").Value = 14.11915767939937
$ws.Range("I2").Value = 22.63773957893154
$ws.Range("K2").Value = 8.015197667241976
$ws.Range("L2").Value = 11.14995757950668
$ws.Range("M2").Value = 13.88875234472146
$ws.Range("N2").Value = 19.0018222748998
$ws.Range("O2").Value = 21.33067805658498

$ws.Range("B3").Value = 10.17610375758471
$ws.Range("C3").Value = 11.0234324933295
$ws.Range("D3").Value = 4.976662779274602
$ws.Range("F3").Value = 23.72176962471077
$ws.Range("G3").Value = 27.74904239640311
$ws.Range("H3").Value = 14.15864578449306
$ws.Range("I3").Value = 22.71647769349634
$ws.Range("K3").Value = 7.76763645695434
$ws.Range("L3").Value = 11.15760807130841
$ws.Range("M3").Value = 13.84983886516642
$ws.Range("N3").Value = 19.05354217578556
$ws.Range("O3").Value = 21.39423429492501

$ws.Range("B4").Value = 10.01707808663902
$ws.Range("C4").Value = 11.02739715878341
$ws.Range("D4").Value = 4.925105955916618
$ws.Range("F4").Value = 23.7484025913675
$ws.Range("G4").Value = 27.7939331711197
$ws.Range("H4").Value = 14.18475552767399
$ws.Range("I4").Value = 22.76844576150417
$ws.Range("K4").Value = 7.610041657735128
$ws.Range("L4").Value = 11.16401047318503
$ws.Range("M4").Value = 13.82787463194401
$ws.Range("N4").Value = 19.08684618217409
$ws.Range("O4").Value = 21.43707274557335

$ws.Range("B5").Value = 9.951918128378665
$ws.Range("C5").Value = 11.02913461345195
$ws.Range("D5").Value = 4.903743219085968
$ws.Range("F5").Value = 23.76062624234852
$ws.Range("G5").Value = 27.81411837815284
$ws.Range("H5").Value = 14.19586450972189
$ws.Range("I5").Value = 22.79053425402713
$ws.Range("K5").Value = 7.544472877578694
$ws.Range("L5").Value = 11.16704888374341
$ws.Range("M5").Value = 13.81941552763858
$ws.Range("N5").Value = 19.10080808498684
$ws.Range("O5").Value = 21.45548833421273

$ws.Range("B6").Value = 9.941079761868583
$ws.Range("C6").Value = 11.02943048810856
$ws.Range("D6").Value = 4.90017501759871
$ws.Range("F6").Value = 23.76273870335944
$ws.Range("G6").Value = 27.81758422996204
$ws.Range("H6").Value = 14.19773748778642
$ws.Range("I6").Value = 22.79425706357963
$ws.Range("K6").Value = 7.533505653301252
$ws.Range("L6").Value = 11.16757935988889
$ws.Range("M6").Value = 13.81804076669274
$ws.Range("N6").Value = 19.10315004933569
$ws.Range("O6").Value = 21.45860410406119

$ws.Range("B7").Value = 10.01620062738
$ws.Range("C7").Value = 11.02742009682203
$ws.Range("D7").Value = 4.924819262423945
$ws.Range("F7").Value = 23.74856189651201
$ws.Range("G7").Value = 27.79419774251675
$ws.Range("H7").Value = 14.18490344764658
$ws.Range("I7").Value = 22.76873996546959
$ws.Range("K7").Value = 7.609162749524967
$ws.Range("L7").Value = 11.16404971092572
$ws.Range("M7").Value = 13.82775855117577
$ws.Range("N7").Value = 19.08703289577613
$ws.Range("O7").Value = 21.43731722431646

$ws.Range("B8").Value = 10.34431240906105
$ws.Range("C8").Value = 11.01961689554282
$ws.Range("D8").Value = 5.030436576043609
$ws.Range("F8").Value = 23.69803773541423
$ws.Range("G8").Value = 27.70761995169143
$ws.Range("H8").Value = 14.13238654588406
$ws.Range("I8").Value = 22.66413672091924
$ws.Range("K8").Value = 7.931032295829155
$ws.Range("L8").Value = 11.15224202390933
$ws.Range("M8").Value = 13.8749386637532
$ws.Range("N8").Value = 19.01933468879401
$ws.Range("O8").Value = 21.35180013398175

$ws.Range("B9").Value = 10.96953659966854
$ws.Range("C9").Value = 11.00812308382576
$ws.Range("D9").Value = 5.225488555247556
$ws.Range("F9").Value = 23.64227264086048
$ws.Range("G9").Value = 27.59786913675904
$ws.Range("H9").Value = 14.04417399720477
$ws.Range("I9").Value = 22.48774579749482
$ws.Range("K9").Value = 8.515570880436316
$ws.Range("L9").Value = 11.14258331049018
$ws.Range("M9").Value = 13.98246156083399
$ws.Range("N9").Value = 18.89881109425016
$ws.Range("O9").Value = 21.21439704215662

$ws.Range("B10").Value = 11.41285338091045
$ws.Range("C10").Value = 11.00197264171389
$ws.Range("D10").Value = 5.360773205853873
$ws.Range("F10").Value = 23.62772410276612
$ws.Range("G10").Value = 27.55404133836499
$ws.Range("H10").Value = 13.98834794199663
$ws.Range("I10").Value = 22.37565870807695
$ws.Range("K10").Value = 8.9139089754068
$ws.Range("L10").Value = 11.14366650412829
$ws.Range("M10").Value = 14.07017796229866
$ws.Range("N10").Value = 18.8176506567581
$ws.Range("O10").Value = 21.13195230274174

$ws.Range("B11").Value = 11.61011439117402
$ws.Range("C11").Value = 10.99966691788597
$ws.Range("D11").Value = 5.420446682212037
$ws.Range("F11").Value = 23.62683615450598
$ws.Range("G11").Value = 27.54212545672234
$ws.Range("H11").Value = 13.96489726648852
$ws.Range("I11").Value = 22.32846690572872
$ws.Range("K11").Value = 9.087895482982928
$ws.Range("L11").Value = 11.14592275299982
$ws.Range("M11").Value = 14.11187722638178
$ws.Range("N11").Value = 18.78231829884551
$ws.Range("O11").Value = 21.09847110986169

$ws.Range("B12").Value = 11.68410845526228
$ws.Range("C12").Value = 10.99886407581584
$ws.Range("D12").Value = 5.442764715882944
$ws.Range("F12").Value = 23.62732237775264
$ws.Range("G12").Value = 27.53876789789021
$ws.Range("H12").Value = 13.95629643609023
$ws.Range("I12").Value = 22.31114256588278
$ws.Range("K12").Value = 9.152708077337529
$ws.Range("L12").Value = 11.14702942862953
$ws.Range("M12").Value = 14.12791710922651
$ws.Range("N12").Value = 18.7691661215663
$ws.Range("O12").Value = 21.08637169456179

$ws.Range("B13").Value = 11.66820490400738
$ws.Range("C13").Value = 10.9990338639838
$ws.Range("D13").Value = 5.437970696102636
$ws.Range("F13").Value = 23.62718111172275
$ws.Range("G13").Value = 27.53943964053416
$ws.Range("H13").Value = 13.95813635482759
$ws.Range("I13").Value = 22.31484937898361
$ws.Range("K13").Value = 9.138797718398314
$ws.Range("L13").Value = 11.14677988801595
$ws.Range("M13").Value = 14.12445169468887
$ws.Range("N13").Value = 18.77198857828464
$ws.Range("O13").Value = 21.08895175128842

$ws.Range("B14").Value = 11.61621642961155
$ws.Range("C14").Value = 10.99959946173897
$ws.Range("D14").Value = 5.422288449456607
$ws.Range("F14").Value = 23.62685968364789
$ws.Range("G14").Value = 27.54182608172902
$ws.Range("H14").Value = 13.96418407214309
$ws.Range("I14").Value = 22.32703067493373
$ws.Range("K14").Value = 9.093249326363704
$ws.Range("L14").Value = 11.14600875167008
$ws.Range("M14").Value = 14.11319189443213
$ws.Range("N14").Value = 18.78123170871906
$ws.Range("O14").Value = 21.09746407253525

$ws.Range("B15").Value = 11.58427825150551
$ws.Range("C15").Value = 10.99995504554296
$ws.Range("D15").Value = 5.412645986792483
$ws.Range("F15").Value = 23.62676985441563
$ws.Range("G15").Value = 27.54343824806222
$ws.Range("H15").Value = 13.96792485674637
$ws.Range("I15").Value = 22.33456320242037
$ws.Range("K15").Value = 9.065209028483068
$ws.Range("L15").Value = 11.14556922218484
$ws.Range("M15").Value = 14.10632713004261
$ws.Range("N15").Value = 18.78692298257917
$ws.Range("O15").Value = 21.10275355989312

$ws.Range("B16").Value = 11.39986706313092
$ws.Range("C16").Value = 11.00213318803152
$ws.Range("D16").Value = 5.356834910562728
$ws.Range("F16").Value = 23.62789732138603
$ws.Range("G16").Value = 27.5549816080526
$ws.Range("H16").Value = 13.98991961421199
$ws.Range("I16").Value = 22.37881922423782
$ws.Range("K16").Value = 8.902390259306127
$ws.Range("L16").Value = 11.14355443033537
$ws.Range("M16").Value = 14.06748820152447
$ws.Range("N16").Value = 18.81999158355239
$ws.Range("O16").Value = 21.13422139301783

$ws.Range("B17").Value = 11.28555605037154
$ws.Range("C17").Value = 11.00359511374657
$ws.Range("D17").Value = 5.322110562718471
$ws.Range("F17").Value = 23.63005572861346
$ws.Range("G17").Value = 27.56411874199881
$ws.Range("H17").Value = 14.00391064466341
$ws.Range("I17").Value = 22.40694150707924
$ws.Range("K17").Value = 8.800631514196047
$ws.Range("L17").Value = 11.1427692063166
$ws.Range("M17").Value = 14.04411558163724
$ws.Range("N17").Value = 18.84068414649835
$ws.Range("O17").Value = 21.15455683454886

$ws.Range("B18").Value = 11.21939652699356
$ws.Range("C18").Value = 11.00448232666789
$ws.Range("D18").Value = 5.301962775399431
$ws.Range("F18").Value = 23.6318366732996
$ws.Range("G18").Value = 27.57012918824928
$ws.Range("H18").Value = 14.01214099064235
$ws.Range("I18").Value = 22.42347409964668
$ws.Range("K18").Value = 8.741424813427702
$ws.Range("L18").Value = 11.14248364123992
$ws.Range("M18").Value = 14.0308419406717
$ws.Range("N18").Value = 18.85273547955846
$ws.Range("O18").Value = 21.16663190388078

$ws.Range("B19").Value = 11.19692773803556
$ws.Range("C19").Value = 11.00479069742082
$ws.Range("D19").Value = 5.295111305539877
$ws.Range("F19").Value = 23.63253236532818
$ws.Range("G19").Value = 27.5722938448038
$ws.Range("H19").Value = 14.01495909504075
$ws.Range("I19").Value = 22.42913313088276
$ws.Range("K19").Value = 8.721263143962693
$ws.Range("L19").Value = 11.14241551506084
$ws.Range("M19").Value = 14.02637712653081
$ws.Range("N19").Value = 18.85684156205523
$ws.Range("O19").Value = 21.17078532986238

$ws.Range("B20").Value = 11.29776770591173
$ws.Range("C20").Value = 11.00343469551221
$ws.Range("D20").Value = 5.325825244715102
$ws.Range("F20").Value = 23.6297701370367
$ws.Range("G20").Value = 27.56306792646977
$ws.Range("H20").Value = 14.00240232971439
$ws.Range("I20").Value = 22.40391085157112
$ws.Range("K20").Value = 8.811534312555029
$ws.Range("L20").Value = 11.14283561567335
$ws.Range("M20").Value = 14.04658614127294
$ws.Range("N20").Value = 18.83846591999863
$ws.Range("O20").Value = 21.15235289725275

$ws.Range("B21").Value = 11.63150635761084
$ws.Range("C21").Value = 10.99943142830746
$ws.Range("D21").Value = 5.426902357550204
$ws.Range("F21").Value = 23.62693178828192
$ws.Range("G21").Value = 27.54109378136916
$ws.Range("H21").Value = 13.96240012918649
$ws.Range("I21").Value = 22.32343790893048
$ws.Range("K21").Value = 9.106657347932522
$ws.Range("L21").Value = 11.14622841725183
$ws.Range("M21").Value = 14.1164924793126
$ws.Range("N21").Value = 18.77851061244392
$ws.Range("O21").Value = 21.09494807377632

$ws.Range("B22").Value = 11.8454892963015
$ws.Range("C22").Value = 10.99722451346707
$ws.Range("D22").Value = 5.491331384454536
$ws.Range("F22").Value = 23.62986963747202
$ws.Range("G22").Value = 27.53346303316914
$ws.Range("H22").Value = 13.93788500811423
$ws.Range("I22").Value = 22.27402776509649
$ws.Range("K22").Value = 9.293273944608551
$ws.Range("L22").Value = 11.1499156130502
$ws.Range("M22").Value = 14.16362886004233
$ws.Range("N22").Value = 18.74065150965567
$ws.Range("O22").Value = 21.06080688801846

$ws.Range("B23").Value = 11.73168279164441
$ws.Range("C23").Value = 10.99836508103245
$ws.Range("D23").Value = 5.457096881377648
$ws.Range("F23").Value = 23.62786375623709
$ws.Range("G23").Value = 27.53691965864957
$ws.Range("H23").Value = 13.95082024898229
$ws.Range("I23").Value = 22.30010754262559
$ws.Range("K23").Value = 9.194256511756038
$ws.Range("L23").Value = 11.14781366431538
$ws.Range("M23").Value = 14.13834182297911
$ws.Range("N23").Value = 18.76073667542697
$ws.Range("O23").Value = 21.07871956781268

$ws.Range("B24").Value = 11.29224818336555
$ws.Range("C24").Value = 11.00350707499694
$ws.Range("D24").Value = 5.324146410239332
$ws.Range("F24").Value = 23.6298975706892
$ws.Range("G24").Value = 27.5635406410916
$ws.Range("H24").Value = 14.00308365727859
$ws.Range("I24").Value = 22.40527987487447
$ws.Range("K24").Value = 8.80660734836111
$ws.Range("L24").Value = 11.14280507524171
$ws.Range("M24").Value = 14.04546869101535
$ws.Range("N24").Value = 18.83946829761584
$ws.Range("O24").Value = 21.15334810137309

$ws.Range("B25").Value = 10.80287690250524
$ws.Range("C25").Value = 11.01082760512024
$ws.Range("D25").Value = 5.174084243662575
$ws.Range("F25").Value = 23.65271780051714
$ws.Range("G25").Value = 27.62110855831821
$ws.Range("H25").Value = 14.06645859317687
$ws.Range("I25").Value = 22.53238900892423
$ws.Range("K25").Value = 8.362728240266353
$ws.Range("L25").Value = 11.14375588864181
$ws.Range("M25").Value = 13.95181015152189
$ws.Range("N25").Value = 18.93011352681744
$ws.Range("O25").Value = 21.24832090184001
